# Insert a new data row at row 23 (pushing existing rows 23-93 down to 24-94)
# and populate it with the new Papaya price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 23; Excel shifts rows 23:93 down to 24:94
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new record's values.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44980
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100108
$ws.Range("H23").Value = "Tropicales y subtropicales"
$ws.Range("I23").Value = 100108004
$ws.Range("J23").Value = "Papaya"
$ws.Range("K23").Value = "Cultivar IV Región"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 80
$ws.Range("N23").Value = 40000
$ws.Range("O23").Value = 40000
$ws.Range("P23").Value = 40000
$ws.Range("Q23").Value = "$/caja 15 kilos granel"
$ws.Range("R23").Value = "Provincia del Elquí"
$ws.Range("S23").Value = 2667
$ws.Range("T23").Value = 15
